# Update countries & provincias Spain
# Refreshes the COVID-19 "Pais" data table with newer figures and updates
# the "datos actualizados" timestamp. A handful of countries changed rank
# (their case counts overtook a neighbouring row), so both the country
# name and the numeric columns for those rows are rewritten.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title / timestamp (row 1) ---------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 18 de Septiembre de 2020 a las 08:56"

# --- Simple numeric refreshes (country order unchanged) ---------------
# Israel
$ws.Range("B27").Value = 176933
$ws.Range("C27").Value = 1677
$ws.Range("D27").Value = 129394
$ws.Range("E27").Value = 46370

# Ucrania
$ws.Range("B28").Value = 169472
$ws.Range("C28").Value = 3228
$ws.Range("D28").Value = 75486
$ws.Range("E28").Value = 90518
$ws.Range("G28").Value = 68
$ws.Range("H28").Value = 3468

# Armenia
$ws.Range("B62").Value = 46910
$ws.Range("C62").Value = 239
$ws.Range("D62").Value = 42369
$ws.Range("E62").Value = 3615
$ws.Range("G62").Value = 1
$ws.Range("H62").Value = 926

# Kirguistan
$ws.Range("B64").Value = 45244
$ws.Range("C64").Value = 91
$ws.Range("D64").Value = 41415
$ws.Range("E64").Value = 2766

# El Salvador
$ws.Range("D75").Value = 20403
$ws.Range("E75").Value = 6042
$ws.Range("G75").Value = 3
$ws.Range("H75").Value = 804

# Australia
$ws.Range("D76").Value = 23855
$ws.Range("E76").Value = 2169

# --- Hungria overtakes Madagascar (rows 86/87 swap) --------------------
$ws.Range("A86").Value = "Hungria"
$ws.Range("B86").Value = 16111
$ws.Range("C86").Value = 941
$ws.Range("D86").Value = 4240
$ws.Range("E86").Value = 11202
$ws.Range("G86").Value = 6
$ws.Range("H86").Value = 669

$ws.Range("A87").Value = "Madagascar"
$ws.Range("B87").Value = 15925
$ws.Range("C87").Value = 0
$ws.Range("D87").Value = 14547
$ws.Range("E87").Value = 1162
$ws.Range("G87").Value = 0
$ws.Range("H87").Value = 216

# --- Georgia overtakes Reunion and Mali (rows 141/142/143 shift) -------
$ws.Range("A141").Value = "Georgia"
$ws.Range("B141").Value = 3119
$ws.Range("C141").Value = 182
$ws.Range("D141").Value = 1435
$ws.Range("E141").Value = 1665
$ws.Range("H141").Value = 19

$ws.Range("A142").Value = "Reunion"
$ws.Range("B142").Value = 3099
$ws.Range("C142").Value = 0
$ws.Range("D142").Value = 1794
$ws.Range("E142").Value = 1290
$ws.Range("H142").Value = 15

$ws.Range("A143").Value = "Mali"
$ws.Range("B143").Value = 2966
$ws.Range("C143").Value = 0
$ws.Range("D143").Value = 2311
$ws.Range("E143").Value = 527
$ws.Range("H143").Value = 128

# Butan
$ws.Range("B187").Value = 252
$ws.Range("C187").Value = 6
$ws.Range("D187").Value = 182
$ws.Range("E187").Value = 70

# --- Montserrat overtakes Islas Malvinas (rows 214/215 swap) -----------
$ws.Range("A214").Value = "Montserrat"
$ws.Range("D214").Value = 12
$ws.Range("H214").Value = 1

$ws.Range("A215").Value = "Islas Malvinas"
$ws.Range("D215").Value = 13
$ws.Range("H215").Value = 0
